# RP2040-ProMini BOM update: U4's footprint was corrected from the
# SOIC-8-208mil package to the WSON-8-EP(6x8) package actually used,
# as part of regenerating the CAM, BOM and CPL files.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Footprint cell for U4 (row 8, column C) to reflect the
# correct WSON-8 footprint.
$ws.Range("C8").Value = "WSON-8-EP(6x8)"

# Reflect the author's final cursor/selection position in the sheet.
$ws.Range("C9").Select()
